# The workbook shipped with image paths like "/img/<dish>.png" in column B
# of Sheet1; the re-upload moves the images under a "/weiwei" sub-folder,
# so every path becomes "/weiwei/img/<dish>.png". It also leaves the
# worksheet's active selection on F8 (instead of B8) and widens column B
# so the longer paths are easier to read.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Re-point every image path in B2:B8 to the new "/weiwei" folder.
$replaced = $ws.Range("B2:B8").Replace("/img/", "/weiwei/img/", 2)

# 2. Update the sheet's remembered selection to F8.
$selected = $ws.Range("F8").Select()

# 3. Widen column B to fit the longer paths.
$ws.Columns("B").ColumnWidth = 26.5
